# Applies the recorded cell-value updates across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets
# (profit-tracking data refreshed by the scheduled market-price runner).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(4, "H").Value = 254
$ws.Cells.Item(4, "I").Value = 276.8
$ws.Cells.Item(4, "K").Value = 276.8
$ws.Cells.Item(4, "M").Value = -162.8
$ws.Cells.Item(131, "H").Value = 4020.6956
$ws.Cells.Item(131, "I").Value = 3725.818
$ws.Cells.Item(131, "J").Value = 4291
$ws.Cells.Item(131, "K").Value = 11177.454
$ws.Cells.Item(131, "L").Value = 12873
$ws.Cells.Item(131, "M").Value = -6137.454000000002
$ws.Cells.Item(131, "N").Value = -22953
$ws.Cells.Item(137, "H").Value = 77250.25
$ws.Cells.Item(137, "I").Value = 100666.664
$ws.Cells.Item(137, "J").Value = 63200.4
$ws.Cells.Item(137, "K").Value = 301999.992
$ws.Cells.Item(137, "L").Value = 189601.2
$ws.Cells.Item(137, "M").Value = -299449.992
$ws.Cells.Item(137, "N").Value = -194701.2
$ws.Cells.Item(138, "H").Value = 27312.65
$ws.Cells.Item(138, "I").Value = 1694.3334
$ws.Cells.Item(138, "J").Value = 65740.125
$ws.Cells.Item(138, "K").Value = 5083.0002
$ws.Cells.Item(138, "L").Value = 197220.375
$ws.Cells.Item(138, "M").Value = 56.9997999999996
$ws.Cells.Item(138, "N").Value = -207500.375

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, "H").Value = 27126.5
$ws.Cells.Item(32, "I").Value = 28546.54
$ws.Cells.Item(32, "K").Value = 28546.54
$ws.Cells.Item(32, "M").Value = -28259.54
$ws.Cells.Item(74, "H").Value = 471660.38
$ws.Cells.Item(74, "I").Value = 1200801.4
$ws.Cells.Item(74, "K").Value = 1200801.4
$ws.Cells.Item(74, "M").Value = -1199927.4
$ws.Cells.Item(77, "H").Value = 471660.38
$ws.Cells.Item(77, "I").Value = 1200801.4
$ws.Cells.Item(77, "K").Value = 6004007
$ws.Cells.Item(77, "M").Value = -5999639
$ws.Cells.Item(97, "H").Value = 683.8182
$ws.Cells.Item(97, "J").Value = 1196
$ws.Cells.Item(97, "L").Value = 1196
$ws.Cells.Item(97, "N").Value = -2188
$ws.Cells.Item(115, "H").Value = 34975
$ws.Cells.Item(115, "J").Value = 34975
$ws.Cells.Item(115, "L").Value = 34975
$ws.Cells.Item(115, "N").Value = -38109
$ws.Cells.Item(118, "H").Value = 100000
$ws.Cells.Item(118, "J").Value = 100000
$ws.Cells.Item(118, "L").Value = 100000
$ws.Cells.Item(118, "N").Value = -103314
$ws.Cells.Item(122, "H").Value = 1478.8518
$ws.Cells.Item(122, "I").Value = 1260
$ws.Cells.Item(122, "K").Value = 3780
$ws.Cells.Item(122, "M").Value = -1330
$ws.Cells.Item(132, "H").Value = 2701
$ws.Cells.Item(132, "I").Value = 2248.0908
$ws.Cells.Item(132, "J").Value = 3697.4
$ws.Cells.Item(132, "K").Value = 6744.2724
$ws.Cells.Item(132, "L").Value = 11092.2
$ws.Cells.Item(132, "M").Value = -4214.2724
$ws.Cells.Item(132, "N").Value = -16152.2

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(20, "H").Value = 2419.3333
$ws.Cells.Item(20, "I").Value = 1833.3334
$ws.Cells.Item(20, "J").Value = 3005.3333
$ws.Cells.Item(20, "K").Value = 1833.3334
$ws.Cells.Item(20, "L").Value = 3005.3333
$ws.Cells.Item(20, "M").Value = -1586.3334
$ws.Cells.Item(20, "N").Value = -3499.3333
$ws.Cells.Item(115, "H").Value = 30000
$ws.Cells.Item(115, "J").Value = 30000
$ws.Cells.Item(115, "L").Value = 30000
$ws.Cells.Item(115, "N").Value = -33134

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(7, "H").Value = 116
$ws.Cells.Item(7, "I").Value = 80
$ws.Cells.Item(7, "J").Value = 143
$ws.Cells.Item(7, "K").Value = 80
$ws.Cells.Item(7, "L").Value = 143
$ws.Cells.Item(7, "M").Value = 33
$ws.Cells.Item(7, "N").Value = -369
$ws.Cells.Item(11, "H").Value = 5000
$ws.Cells.Item(11, "J").Value = 5000
$ws.Cells.Item(11, "L").Value = 5000
$ws.Cells.Item(11, "N").Value = -5280
$ws.Cells.Item(31, "H").Value = 5883357
$ws.Cells.Item(31, "J").Value = 1554.5
$ws.Cells.Item(31, "L").Value = 1554.5
$ws.Cells.Item(31, "N").Value = -2144.5
$ws.Cells.Item(34, "H").Value = 5883357
$ws.Cells.Item(34, "J").Value = 1554.5
$ws.Cells.Item(34, "L").Value = 1554.5
$ws.Cells.Item(34, "N").Value = -1958.5
$ws.Cells.Item(94, "H").Value = 1559.381
$ws.Cells.Item(94, "I").Value = 1229.3636
$ws.Cells.Item(94, "K").Value = 1229.3636
$ws.Cells.Item(94, "M").Value = -778.3635999999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(4, "H").Value = 100804970
$ws.Cells.Item(4, "I").Value = 50153868
$ws.Cells.Item(4, "J").Value = 333800000
$ws.Cells.Item(4, "K").Value = 150461604
$ws.Cells.Item(4, "L").Value = 1001400000
$ws.Cells.Item(4, "M").Value = -150461492
$ws.Cells.Item(4, "N").Value = -1001400224
$ws.Cells.Item(37, "H").Value = 41972.555
$ws.Cells.Item(37, "J").Value = 41972.555
$ws.Cells.Item(37, "L").Value = 125917.665
$ws.Cells.Item(37, "N").Value = -126141.665
$ws.Cells.Item(80, "H").Value = 4999.2383
$ws.Cells.Item(80, "J").Value = 4999.2383
$ws.Cells.Item(80, "L").Value = 14997.7149
$ws.Cells.Item(80, "N").Value = -16869.7149
$ws.Cells.Item(83, "H").Value = 4999.2383
$ws.Cells.Item(83, "J").Value = 4999.2383
$ws.Cells.Item(83, "L").Value = 44993.1447
$ws.Cells.Item(83, "N").Value = -54353.1447
$ws.Cells.Item(92, "H").Value = 1530.3846
$ws.Cells.Item(92, "I").Value = 1534.6
$ws.Cells.Item(92, "J").Value = 1516.3334
$ws.Cells.Item(92, "K").Value = 4603.799999999999
$ws.Cells.Item(92, "L").Value = 4549.0002
$ws.Cells.Item(92, "M").Value = -3355.799999999999
$ws.Cells.Item(92, "N").Value = -7045.0002
$ws.Cells.Item(109, "H").Value = 890.1
$ws.Cells.Item(127, "H").Value = 4662.4443
$ws.Cells.Item(127, "I").Value = 1030
$ws.Cells.Item(127, "J").Value = 5116.5
$ws.Cells.Item(127, "K").Value = 3090
$ws.Cells.Item(127, "L").Value = 15349.5
$ws.Cells.Item(127, "M").Value = 1870
$ws.Cells.Item(127, "N").Value = -25269.5
$ws.Cells.Item(140, "H").Value = 3258.75
$ws.Cells.Item(140, "I").Value = 3258.75
$ws.Cells.Item(140, "K").Value = 9776.25
$ws.Cells.Item(140, "M").Value = -4596.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(70, "H").Value = 6064.8887
$ws.Cells.Item(70, "I").Value = 7184.25
$ws.Cells.Item(70, "K").Value = 7184.25
$ws.Cells.Item(70, "M").Value = -6914.25
$ws.Cells.Item(73, "H").Value = 6064.8887
$ws.Cells.Item(73, "I").Value = 7184.25
$ws.Cells.Item(73, "K").Value = 7184.25
$ws.Cells.Item(73, "M").Value = -6248.25
$ws.Cells.Item(97, "H").Value = 936.4138
$ws.Cells.Item(97, "I").Value = 963.3889
$ws.Cells.Item(97, "K").Value = 963.3889
$ws.Cells.Item(97, "M").Value = -467.3889
$ws.Cells.Item(102, "H").Value = 16836.688
$ws.Cells.Item(102, "I").Value = 19467.111
$ws.Cells.Item(102, "K").Value = 19467.111
$ws.Cells.Item(102, "M").Value = -17845.111
$ws.Cells.Item(122, "H").Value = 3044.074
$ws.Cells.Item(122, "I").Value = 2840.476
$ws.Cells.Item(122, "K").Value = 8521.428
$ws.Cells.Item(122, "M").Value = -6071.428
$ws.Cells.Item(132, "H").Value = 2241.5334
$ws.Cells.Item(132, "I").Value = 2100.818
$ws.Cells.Item(132, "J").Value = 2628.5
$ws.Cells.Item(132, "K").Value = 6302.454000000001
$ws.Cells.Item(132, "L").Value = 7885.5
$ws.Cells.Item(132, "M").Value = -3772.454000000001
$ws.Cells.Item(132, "N").Value = -12945.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(2, "H").Value = 15000000
$ws.Cells.Item(2, "J").Value = 15000000
$ws.Cells.Item(2, "L").Value = 15000000
$ws.Cells.Item(2, "N").Value = -15000224
$ws.Cells.Item(40, "H").Value = 4656.625
$ws.Cells.Item(40, "I").Value = 4763
$ws.Cells.Item(40, "J").Value = 3912
$ws.Cells.Item(40, "K").Value = 4763
$ws.Cells.Item(40, "L").Value = 3912
$ws.Cells.Item(40, "M").Value = -4627
$ws.Cells.Item(40, "N").Value = -4184
$ws.Cells.Item(43, "H").Value = 13636.363
$ws.Cells.Item(43, "I").Value = 0
$ws.Cells.Item(43, "J").Value = 13636.363
$ws.Cells.Item(43, "K").Value = 0
$ws.Cells.Item(43, "L").Value = 13636.363
$ws.Cells.Item(43, "M").ClearContents()
$ws.Cells.Item(43, "N").Value = -14022.363
$ws.Cells.Item(136, "H").Value = 3328.625
$ws.Cells.Item(136, "I").Value = 3393.1
$ws.Cells.Item(136, "J").Value = 3221.1667
$ws.Cells.Item(136, "K").Value = 10179.3
$ws.Cells.Item(136, "L").Value = 9663.500100000001
$ws.Cells.Item(136, "M").Value = -7629.299999999999
$ws.Cells.Item(136, "N").Value = -14763.5001

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(2, "H").Value = 460.5
$ws.Cells.Item(2, "I").Value = 460.5
$ws.Cells.Item(2, "K").Value = 460.5
$ws.Cells.Item(2, "M").Value = -348.5
$ws.Cells.Item(95, "H").Value = 57797
$ws.Cells.Item(95, "J").Value = 57797
$ws.Cells.Item(95, "L").Value = 57797
$ws.Cells.Item(95, "N").Value = -63289
$ws.Cells.Item(107, "H").Value = 913.5625
$ws.Cells.Item(107, "I").Value = 1399
$ws.Cells.Item(107, "K").Value = 4197
$ws.Cells.Item(107, "M").Value = -2277
$ws.Cells.Item(136, "H").Value = 26631.72
$ws.Cells.Item(136, "I").Value = 29865.592
$ws.Cells.Item(136, "K").Value = 89596.776
$ws.Cells.Item(136, "M").Value = -87046.776
